# Commit: commited on with latest changes on tickets
#
# Summary of the change being reproduced:
#  - A new worksheet "Ticket" is inserted right after "Project" (and before
#    "Message"), becomes the active sheet/tab.
#  - The Ticket sheet gets a small 2-row table: Title/Description/Client
#    headers with TicketVM / To add ticket / 123 data, plus a couple of
#    custom column widths.
#  - On the "Message" sheet, the value that used to read "TestSample15" is
#    corrected to "TestSample".
#  - Leftover cursor/selection state changes on the "Project" and "Message"
#    sheets (captured by the original author's workbook save).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Ticket" worksheet directly after "Project" and fill
#    in its little lookup table. The shared-string table records brand
#    new text in first-use order, so write these cells (TicketVM / To add
#    ticket / Client) before the "TestSample" correction below to match
#    the saved file's <sst> ordering.
# ---------------------------------------------------------------------
$project = $wb.Worksheets.Item("Project")
$ticket = $wb.Worksheets.Add($null, $project)
$ticket.Name = "Ticket"

$ticket.Range("A1").Value = "Title"
$ticket.Range("B1").Value = "Description"
$ticket.Range("A2").Value = "TicketVM"
$ticket.Range("B2").Value = "To add ticket"
$ticket.Range("C2").Value = 123

$ticket.Columns.Item(1).ColumnWidth = 11.666666666666666
$ticket.Columns.Item(2).ColumnWidth = 12.333333333333332

# ---------------------------------------------------------------------
# 2. Update the "Message" sheet cell content (TestSample15 -> TestSample).
# ---------------------------------------------------------------------
$message = $wb.Worksheets.Item("Message")
$message.Range("B2").Value = "TestSample"

# ---------------------------------------------------------------------
# 3. Last new shared string introduced in the saved file is "Client".
# ---------------------------------------------------------------------
$ticket.Range("C1").Value = "Client"

# ---------------------------------------------------------------------
# 4. Leftover cursor/selection state, as captured by the original
#    author's workbook save: Project -> G16, Message -> B2, and the
#    newly added Ticket sheet ends up active/selected at G5.
# ---------------------------------------------------------------------
$project.Range("G16").Select()
$message.Range("B2").Select()
$ticket.Range("G5").Select()
